$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Channel column (C) is repurposed: every data row (2-55) now refers to the
# "Canteen_OTH" channel instead of the previous mix of
# "Canteen_TRAD" / "Canteen" / "Supermarket_CAP" values.
$ws.Range("C2:C55").Value = "Canteen_OTH"

# Register the new (auto-numbered) hidden filter-database name that Excel
# creates when the existing AutoFilter range is touched/re-saved.
$ws.Names.Add("_xlnm._FilterDatabase_0_0", "=Canteen!`$A`$1:`$AL`$1")

# Restore/update the view state: scroll back to the top of the frozen pane
# and leave the edited range selected.
$ws.Range("C2:C55").Select()
